$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C3 value from 1 to 0
$ws.Range("C3").Value = 0

# Update selection to C3 (reflects the <selection activeCell="C3" sqref="C3"/> in the saved file)
$ws.Range("C3").Select()
